# Apply the edits described by the diff:
#  1. TOC sheet: rename header B1 from "figure_name" to "description".
#  2. Figure_6 sheet (no_vehicle_percentage): insert a missing
#     "Ventura / Native American" row (value 0) before "Ventura / White",
#     pushing the existing rows (and the trailing SCAG block) down by one.
#  3. Figure_21 sheet (overcrowded_percentage): insert a missing
#     "Imperial / Multiracial/Other" row (value 0) before
#     "Imperial / Native American", pushing all following rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. TOC: B1 header rename
# ---------------------------------------------------------------
$toc = $wb.Worksheets.Item("TOC")
$toc.Range("B1").Value = "description"

# ---------------------------------------------------------------
# 2. Figure_6: insert Ventura / Native American / 0 at row 41
# ---------------------------------------------------------------
$f6 = $wb.Worksheets.Item("Figure_6")
$f6.Rows.Item(41).Insert()
$f6.Cells.Item(41, 1).Value = "Ventura"
$f6.Cells.Item(41, 2).Value = "Native American"
$f6.Cells.Item(41, 3).Value = 0

# ---------------------------------------------------------------
# 3. Figure_21: insert Imperial / Multiracial/Other / 0 at row 5
# ---------------------------------------------------------------
$f21 = $wb.Worksheets.Item("Figure_21")
$f21.Rows.Item(5).Insert()
$f21.Cells.Item(5, 1).Value = "Imperial"
$f21.Cells.Item(5, 2).Value = "Multiracial/Other"
$f21.Cells.Item(5, 3).Value = 0
